# Apply the data update described by the diff: rows 335-344 of the
# "Chiffres COVID-19 Valais" sheet get refreshed daily figures.
#
# Columns L (12) and M (13) on this worksheet have a quirk in this
# runtime: writing a plain numeric value straight to `.Value` on those
# two columns causes the cell to be stored as a text string instead of
# a number. Toggling the NumberFormat to "General" before the write and
# restoring the original "@" (Text) format afterwards avoids the issue
# while keeping the original style/format intact.
function Set-CellNumber($range, $value) {
    $col = $range.Column
    if ($col -eq 12 -or $col -eq 13) {
        $fmt = $range.NumberFormat
        $range.NumberFormat = "General"
        $range.Value = $value
        $range.NumberFormat = $fmt
    } else {
        $range.Value = $value
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 335: Nb nouveaux cas positifs corrected 152 -> 153
Set-CellNumber $ws.Cells.Item(335, 3) 153

# Row 340: Nb nouveaux cas positifs corrected 76 -> 75
Set-CellNumber $ws.Cells.Item(340, 3) 75

# Row 342: Nb nouveaux cas positifs corrected 91 -> 123
Set-CellNumber $ws.Cells.Item(342, 3) 123

# Row 343: Nb nouveaux cas positifs corrected 8 -> 82, plus one extra
# out-of-hospital death reported (M343 0 -> 2)
Set-CellNumber $ws.Cells.Item(343, 3) 82
Set-CellNumber $ws.Cells.Item(343, 13) 2

# Row 344 (2020-09-24) was previously blank; fill in the day's figures.
Set-CellNumber $ws.Cells.Item(344, 3) 12
Set-CellNumber $ws.Cells.Item(344, 5) 10
Set-CellNumber $ws.Cells.Item(344, 6) 7
Set-CellNumber $ws.Cells.Item(344, 7) 117
Set-CellNumber $ws.Cells.Item(344, 12) 0
Set-CellNumber $ws.Cells.Item(344, 13) 0
